$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 25 - this shifts existing rows 25-73 down to 26-74,
# matching the new dimension A1:R74, and carries the date-column style (s=2)
# on column D down into the new row automatically (Excel's native InsertShift behavior).
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the weekly data point.
# Columns A,B,C,E,F,G,H,I,N,Q,R keep the same values the (now shifted) row below had;
# D,J,K,L,M,O,P carry the new week's figures.
$ws.Range("A25").Value = 11
$ws.Range("B25").Value = "Vega Monumental Concepción"
$ws.Range("C25").Value = "Bíobío"
$ws.Range("D25").Value = 45238
$ws.Range("E25").Value = 8
$ws.Range("F25").Value = 100112026
$ws.Range("G25").Value = "Haba"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 50
$ws.Range("K25").Value = 12000
$ws.Range("L25").Value = 12000
$ws.Range("M25").Value = 12000
$ws.Range("N25").Value = "$/saco 25 kilos"
$ws.Range("O25").Value = "Región del Maule"
$ws.Range("P25").Value = 480
$ws.Range("Q25").Value = 25
$ws.Range("R25").Value = "Hortaliza"
